$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving default (unstyled) cell style
function Set-TextValue($ws, $addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# Row 2
Set-TextValue $ws "D2" "43.772.31"
Set-TextValue $ws "E2" "  +4.87%  "

# Row 3
Set-TextValue $ws "D3" "2.273.20"
Set-TextValue $ws "E3" "  +1.84%  "

# Row 4
Set-TextValue $ws "E4" "  +0.14%  "

# Row 5
Set-TextValue $ws "D5" "231.15"
Set-TextValue $ws "E5" "  -0.04%  "

# Row 6
Set-TextValue $ws "E6" "  +0.64%  "

# Row 7
Set-TextValue $ws "D7" "63.91"
Set-TextValue $ws "E7" "  +6.30%  "

# Row 8
Set-TextValue $ws "E8" "  +0.07%  "

# Row 9
Set-TextValue $ws "D9" "0.426"
Set-TextValue $ws "E9" "  +5.79%  "

# Row 10
Set-TextValue $ws "E10" "  +16.38%  "

# Row 11
Set-TextValue $ws "D11" "57.27"
Set-TextValue $ws "E11" "  -1.15%  "

# Row 12
Set-TextValue $ws "D12" "25.93"
Set-TextValue $ws "E12" "  +15.36%  "

# Row 13
Set-TextValue $ws "E13" "  -0.03%  "

# Row 14
Set-TextValue $ws "D14" "2.610.87"
Set-TextValue $ws "E14" "  +1.89%  "

# Row 15
Set-TextValue $ws "E15" "  +1.29%  "

# Row 16
Set-TextValue $ws "E16" "  +4.19%  "

# Row 17
Set-TextValue $ws "D17" "0.817"
Set-TextValue $ws "E17" "  +2.21%  "

# Row 18
Set-TextValue $ws "D18" "2.280.62"
Set-TextValue $ws "E18" "  +1.50%  "

# Row 19
Set-TextValue $ws "D19" "43.686.35"

# Row 20
Set-TextValue $ws "D20" "0.0000101"
Set-TextValue $ws "E20" "  +10.71%  "

# Row 21
Set-TextValue $ws "D21" "73.19"
Set-TextValue $ws "E21" "  +1.03%  "

# Row 22
Set-TextValue $ws "D22" "6.08"
Set-TextValue $ws "E22" "  -0.86%  "

# Row 23
Set-TextValue $ws "D23" "248.67"
Set-TextValue $ws "E23" "  +0.41%  "

# Row 24
Set-TextValue $ws "E24" "  +0.16%  "

# Row 25
Set-TextValue $ws "D25" "2.49"
Set-TextValue $ws "E25" "  +5.50%  "

# Row 26
Set-TextValue $ws "D26" "2.26"
Set-TextValue $ws "E26" "  -2.04%  "

# Row 27
Set-TextValue $ws "D27" "9.82"
Set-TextValue $ws "E27" "  +0.37%  "

# Row 28
Set-TextValue $ws "D28" "171.62"
Set-TextValue $ws "E28" "  +1.35%  "

# Row 29
Set-TextValue $ws "B29" "EthereumClassic"
Set-TextValue $ws "C29" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws "D29" "20.87"
Set-TextValue $ws "E29" "  +4.89%  "

# Row 30
Set-TextValue $ws "B30" "Kaspa"
Set-TextValue $ws "C30" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws "D30" "0.137"
Set-TextValue $ws "E30" "  -2.58%  "

# Row 31
Set-TextValue $ws "D31" "1.44"
Set-TextValue $ws "E31" "  +2.06%  "

# Row 32
Set-TextValue $ws "D32" "2.79"
Set-TextValue $ws "E32" "  +8.60%  "

# Row 33
Set-TextValue $ws "D33" "0.122"
Set-TextValue $ws "E33" "  -0.11%  "

# Row 34
Set-TextValue $ws "D34" "0.0688"
Set-TextValue $ws "E34" "  +5.16%  "

# Row 35
Set-TextValue $ws "D35" "5.09"
Set-TextValue $ws "E35" "  +2.03%  "

# Row 36
Set-TextValue $ws "D36" "4.69"
Set-TextValue $ws "E36" "  +0.19%  "

# Row 37
Set-TextValue $ws "B37" "RenderToken"
Set-TextValue $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws "D37" "3.81"
Set-TextValue $ws "E37" "  +6.14%  "

# Row 38
Set-TextValue $ws "B38" "THORChain"
Set-TextValue $ws "C38" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws "D38" "6.74"
Set-TextValue $ws "E38" "  +2.94%  "

# Row 39
Set-TextValue $ws "E39" "  -3.12%  "

# Row 40
Set-TextValue $ws "E40" "  +3.66%  "

# Row 41
Set-TextValue $ws "E41" "  +0.24%  "

# Row 42
Set-TextValue $ws "E42" "  -2.86%  "

# Row 43
Set-TextValue $ws "D43" "10.56"
Set-TextValue $ws "E43" "  +19.92%  "

# Row 44
Set-TextValue $ws "B44" "InjectiveProtocol"
Set-TextValue $ws "C44" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws "D44" "17.21"
Set-TextValue $ws "E44" "  +4.02%  "

# Row 45
Set-TextValue $ws "B45" "Cronos"
Set-TextValue $ws "C45" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D45" "0.0962"
Set-TextValue $ws "E45" "  -0.02%  "

# Row 46
Set-TextValue $ws "B46" "TrustWalletToken"
Set-TextValue $ws "C46" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws "D46" "1.20"
Set-TextValue $ws "E46" "  -0.94%  "

# Row 47
Set-TextValue $ws "B47" "FTXToken"
Set-TextValue $ws "C47" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws "D47" "4.42"
Set-TextValue $ws "E47" "  +1.14%  "

# Row 48
Set-TextValue $ws "B48" "Aave"
Set-TextValue $ws "C48" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws "D48" "97.19"
Set-TextValue $ws "E48" "  -1.79%  "

# Row 49
Set-TextValue $ws "D49" "1.475.44"
Set-TextValue $ws "E49" "  -0.14%  "

# Row 50
Set-TextValue $ws "D50" "0.000208"
Set-TextValue $ws "E50" "  -13.65%  "

# Row 51
Set-TextValue $ws "D51" "2.32"
Set-TextValue $ws "E51" "  +3.60%  "
